$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General numeric-looking values must stay text) by setting
# the "@" text number format before assigning values, then resetting the
# style back to Normal so no stray style index is left on the cell.
$textCells = @("D5", "D6", "D8", "D11", "D12", "D13", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

# Apply the new cell values exactly as captured in the target diff.
$ws.Range("D2").Value = "92.980.12"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "3.108.29"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "241.61"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "614.55"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +5.91%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "3.109.50"
$ws.Range("E10").Value = "  +30.17%  "
$ws.Range("D11").Value = "0.754"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "0.201"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "93.153.56"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "34.40"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "3.710.08"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "3.120.42"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "3.77"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "14.80"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "5.78"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "446.79"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").Value = "0.0000202"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "5.77"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "86.93"
$ws.Range("E26").Value = "  +5.69%  "
$ws.Range("D27").Value = "11.71"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "3.293.87"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +10.11%  "
$ws.Range("D31").Value = "0.233"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "0.169"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "9.17"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").Value = "8.07"
$ws.Range("E34").Value = "  +5.94%  "
$ws.Range("D35").Value = "0.159"
$ws.Range("E35").Value = "  -6.69%  "
$ws.Range("D36").Value = "26.09"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  -6.58%  "
$ws.Range("D38").Value = "491.72"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "3.86"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "1.29"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "0.430"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "23.12"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.38"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D46").Value = "163.47"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "0.686"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "44.03"
$ws.Range("E51").Value = "  -0.12%  "

# Reset number format on the forced-text cells back to the default style
# now that the text values are safely committed.
foreach ($addr in $textCells) { $ws.Range($addr).Style = "Normal" }
